# feat: add account name to db
# Adds a new column J containing the account identifier "30991" to every
# data row (1-375) of the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 375
$accountName = "30991"

# Force the new column to be stored as text (not a number) so that the
# value "30991" round-trips as a string, matching the source data.
$rng = $ws.Range("J1:J$lastRow")
$rng.NumberFormat = "@"

for ($r = 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 10).Value = $accountName
}
